$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value from 0.4.0 to 0.7.0
$ws.Range("B3").Value = "0.7.0"

# Delete the entire "Jurisdiction" / "Chile" row (row 11)
$ws.Rows.Item(11).Delete()
